$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "FormVersion"
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = "form_version"
$ws.Range("D41").Value = "FormVersion"

$ws.Range("A41").Style = $ws.Range("A40").Style
$ws.Range("B41").Style = $ws.Range("B9").Style
$ws.Range("C41").Style = $ws.Range("C40").Style
$ws.Range("D41").Style = $ws.Range("D2").Style

$ws.Range("A41").Select()
